# "Do not allow duplicate antibody names"
#
# The validation service re-ran over the Antibodies sheet: a duplicate
# antibody name ("Acme mAb 1") was detected in row 4, which pushed all of
# the already-flagged problems in rows 5-10 down by one row (now rows
# 6-11), and a new "Duplicate antibody name" comment/highlight was added
# for A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Antibodies")

# ---------------------------------------------------------------------
# Helpers: reuse existing, already-present style indices instead of
# creating new ones. A3 carries the pink "problem" highlight (style 2)
# and is never touched by this change, so it is a safe style donor; A2
# is plain/default (style 0) and is likewise untouched.
# ---------------------------------------------------------------------
function Set-Highlight($rangeAddress) {
    $ws.Range("A3").Copy() | Out-Null
    $ws.Range($rangeAddress).PasteSpecial(-4122) | Out-Null
}

function Clear-Highlight($rangeAddress) {
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range($rangeAddress).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 4: "Acme mAb 3" -> "Acme mAb 1" (a duplicate of row 2), now
# highlighted as invalid.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Acme mAb 1"
Set-Highlight "A4"

# ---------------------------------------------------------------------
# Row 5: Host was missing (highlighted blank) - now filled in correctly,
# so the highlight is removed.
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "Homo sapiens"
Clear-Highlight "B5"

# ---------------------------------------------------------------------
# Row 6: Host used to be the unrecognized "Mu musculus" - now blank
# (missing), still highlighted.
# ---------------------------------------------------------------------
$ws.Range("B6").Value = ""
Set-Highlight "B6"

# ---------------------------------------------------------------------
# Row 7: Host used to be the unrecognized "Coronavirus" - now the
# unrecognized "Mu musculus"; stays highlighted either way.
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "Mu musculus"
Set-Highlight "B7"

# ---------------------------------------------------------------------
# Row 8: Host used to be the valid "Homo sapiens" (no highlight) - now
# the unrecognized "Coronavirus", newly highlighted.
# ---------------------------------------------------------------------
$ws.Range("B8").Value = "Coronavirus"
Set-Highlight "B8"

# ---------------------------------------------------------------------
# Row 9: Isotype was missing (highlighted blank) - now filled in
# correctly, so the highlight is removed.
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "IgA2"
Clear-Highlight "C9"

# ---------------------------------------------------------------------
# Row 10: Isotype used to be the unrecognized "Ig" - now blank
# (missing), still highlighted.
# ---------------------------------------------------------------------
$ws.Range("C10").Value = ""
Set-Highlight "C10"

# ---------------------------------------------------------------------
# Row 11: Isotype used to be the valid "IgM" (no highlight) - now the
# unrecognized "Ig", newly highlighted.
# ---------------------------------------------------------------------
$ws.Range("C11").Value = "Ig"
Set-Highlight "C11"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Comments: shift the existing validation comments down to their new
# rows (re-using the comment already sitting in that cell, so its
# original author is preserved), drop the one that no longer applies,
# and add the brand new duplicate-name comment plus the two that now
# land on previously-unflagged cells.
# ---------------------------------------------------------------------
$ws.Range("B5").Comment.Delete()
$ws.Range("C9").Comment.Delete()

$ws.Range("B6").AddComment("Missing required value 'Host'") | Out-Null
$ws.Range("B7").AddComment("'Mu musculus' is not a recognized host") | Out-Null
$ws.Range("C10").AddComment("Missing required value 'Isotype'") | Out-Null

$ws.Range("A4").AddComment("Duplicate antibody name 'Acme mAb 1' is not allowed") | Out-Null
$ws.Range("B8").AddComment("'Coronavirus' is not a recognized host") | Out-Null
$ws.Range("C11").AddComment("'Ig' is not a recognized isotype") | Out-Null
